$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D (shifts existing D:K -> F:M)
$ws.Columns("D:E").Insert()

# Copy number-format/style from column F (the original D, now shifted) onto new D:E
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)

# Rows that have no data in D:K originally (section headers) should not gain D:E cells
$ws.Range("D5:E6").Clear()
$ws.Range("D36:E37").Clear()
$ws.Range("D78:E79").Clear()

$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 1411300
$ws.Range("E8").Value = 1629900
$ws.Range("D9").Value = 984800
$ws.Range("E9").Value = 1151900
$ws.Range("D10").Value = 426500
$ws.Range("E10").Value = 478000
$ws.Range("D12:E12").Value = "NA"
$ws.Range("D13:E13").Value = 0
$ws.Range("D14").Value = -107900
$ws.Range("E14").Value = 133800
$ws.Range("D15").Value = 40300
$ws.Range("E15").Value = 42300
$ws.Range("D17").Value = 1123500
$ws.Range("E17").Value = 1558600
$ws.Range("D18").Value = 287800
$ws.Range("E18").Value = 71300
$ws.Range("D20").Value = -58700
$ws.Range("E20").Value = 24900
$ws.Range("D21").Value = 322700
$ws.Range("E21").Value = 193800
$ws.Range("D22").Value = 59400
$ws.Range("E22").Value = 99100
$ws.Range("D23").Value = 169700
$ws.Range("E23").Value = -2900
$ws.Range("D24").Value = 314700
$ws.Range("E24").Value = 7400
$ws.Range("D25:E25").Value = 0
$ws.Range("D26").Value = -145000
$ws.Range("E26").Value = -10300
$ws.Range("D27").Value = -147300
$ws.Range("E27").Value = -12500
$ws.Range("D28:E28").Value = 0
$ws.Range("D29").Value = 270900
$ws.Range("E29").Value = -5100
$ws.Range("D30:E30").Value = 0
$ws.Range("D31:E31").Value = 0
$ws.Range("D32").Value = 58700
$ws.Range("E32").Value = -24900
$ws.Range("D33").Value = 123600
$ws.Range("E33").Value = -17600
$ws.Range("D34:E34").Value = 0
$ws.Range("D35").Value = 123600
$ws.Range("E35").Value = -17600
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 222900
$ws.Range("E41").Value = 989700
$ws.Range("D42:E42").Value = 0
$ws.Range("D43").Value = 452500
$ws.Range("E43").Value = 462300
$ws.Range("D44").Value = 498900
$ws.Range("E44").Value = 484200
$ws.Range("D45").Value = 68100
$ws.Range("E45").Value = 264100
$ws.Range("D46").Value = 1242400
$ws.Range("E46").Value = 2200300
$ws.Range("D47").Value = 168100
$ws.Range("E47").Value = "NA"
$ws.Range("D48").Value = 1715000
$ws.Range("E48").Value = 1709700
$ws.Range("D49").Value = 7967300
$ws.Range("E49").Value = 8038900
$ws.Range("D50:E50").Value = 0
$ws.Range("D51:E51").Value = 0
$ws.Range("D52").Value = 192000
$ws.Range("E52").Value = 1108600
$ws.Range("D53:E53").Value = 0
$ws.Range("D54").Value = 11284800
$ws.Range("E54").Value = 13057500
$ws.Range("D57").Value = 331600
$ws.Range("E57").Value = 365100
$ws.Range("D58").Value = 3400
$ws.Range("E58").Value = 22100
$ws.Range("D59").Value = 462100
$ws.Range("E59").Value = 404900
$ws.Range("D60").Value = 797100
$ws.Range("E60").Value = 792100
$ws.Range("D61").Value = 6336500
$ws.Range("E61").Value = 7232100
$ws.Range("D62").Value = 994800
$ws.Range("E62").Value = 1972800
$ws.Range("D63:E63").Value = 0
$ws.Range("D64:E64").Value = 0
$ws.Range("D65:E65").Value = 0
$ws.Range("D66").Value = 8138800
$ws.Range("E66").Value = 10007100
$ws.Range("D68:E68").Value = 0
$ws.Range("D69:E69").Value = 0
$ws.Range("D70:E70").Value = 0
$ws.Range("D71:E71").Value = 0
$ws.Range("D72").Value = 210700
$ws.Range("E72").Value = 88000
$ws.Range("D73:E73").Value = 0
$ws.Range("D74:E74").Value = 0
$ws.Range("D75:E75").Value = 0
$ws.Range("D76").Value = 3146000
$ws.Range("E76").Value = 3050400
$ws.Range("D77:E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 123600
$ws.Range("E81").Value = -17600
$ws.Range("D83").Value = 93600
$ws.Range("E83").Value = 97600
$ws.Range("D84:E84").Value = 0
$ws.Range("D85:E85").Value = 0
$ws.Range("D86:E86").Value = 0
$ws.Range("D87:E87").Value = 0
$ws.Range("D88:E88").Value = 0
$ws.Range("D89").Value = 238700
$ws.Range("E89").Value = 128200
$ws.Range("D91").Value = -78800
$ws.Range("E91").Value = -82900
$ws.Range("D92:E92").Value = 0
$ws.Range("D93:E93").Value = 0
$ws.Range("D94").Value = 201500
$ws.Range("E94").Value = -78100
$ws.Range("D96:E96").Value = 0
$ws.Range("D97:E97").Value = 0
$ws.Range("D98:E98").Value = 0
$ws.Range("D99:E99").Value = 0
$ws.Range("D100").Value = -1199600
$ws.Range("E100").Value = 597300
$ws.Range("D101").Value = -1600
$ws.Range("E101").Value = -300
$ws.Range("D102").Value = -761000
$ws.Range("E102").Value = 647100
$ws.Range("F47:J47").Value = "NA"
$ws.Range("I91").Value = -65400
$ws.Range("J91").Value = -125000
